$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Text")
Write-Host $ws.Name
Write-Host $ws.Range("A1").Value
